$d = $word.ActiveDocument

# Build the text for "Nombre d'image" using the actual right single quotation
# mark (U+2019) that the document uses.
$rsquo = [char]0x2019
$nombreImage = "Nombre d" + $rsquo + "image"
$categories = [char]0x0043 + "at" + [char]0x00E9 + "gories"

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($t -eq "Pays") {
        # Strike through the "Pays" bullet (paragraph mark + run).
        $p.Range.Font.StrikeThrough = 1
    }
    elseif ($t -eq $categories) {
        # Highlight the "Cat\u00e9gories" bullet in yellow.
        $p.Range.Font.HighlightColorIndex = 7
    }
    elseif ($t -eq $nombreImage) {
        # Highlight the "Nombre d\u2019image" bullet in yellow.
        $p.Range.Font.HighlightColorIndex = 7
    }
}
